# Applies the "Including predation losses and gains" commit:
#  - renames Feuil1 -> Calcul
#  - adds a new "Output" sheet that pulls together results from Calcul
#  - adds the "Nombre d'armements" / Crowhurst-Acerra-Meyer block (rows 193-231) on Calcul
#  - adds the running totals (B136, B162) and the "Part de captures prives" ratio (D161/D162)

$wb = $excel.ActiveWorkbook

# --- rename the original sheet, add the new one right after it -------------
$calcul = $wb.ActiveSheet
$calcul.Name = "Calcul"

$output = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $calcul)
$output.Name = "Output"

# ============================================================================
# Calcul sheet edits
# ============================================================================

# -- running total under the "guerre d'Amerique" (King's war share) block --
$calcul.Range("B136").Formula = "=SUM(B130:B135)"

# -- "Nombre d'armements (?)" block -----------------------------------------
$calcul.Range("A193").Value = "Crowhurst par Acerra Meyer ""Appendice 1 Marine et Révolution"""

$calcul.Range("B194").Value = "Nombre d’armements (?)"

$calcul.Range("A195").Value = 1793
$calcul.Range("B195").Value = 132

$calcul.Range("A196").Value = 1794

$calcul.Range("A197").Value = 1795

$calcul.Range("A198").Value = 1796
$calcul.Range("B198").Value = 45

$calcul.Range("A199").Value = 1797
$calcul.Range("B199").Value = 145

$calcul.Range("A200").Value = 1798
$calcul.Range("B200").Value = 144

$calcul.Range("A201").Value = 1799
$calcul.Range("B201").Value = 117

$calcul.Range("A202").Value = 1800
$calcul.Range("B202").Value = 87

$calcul.Range("A203").Value = 1801
$calcul.Range("B203").Value = 55

# -- "Dukerque guerre d'Amerique" note + average armament -------------------
$calcul.Range("A205").Value = "Dukerque guerre d’Amérique : 198 armements (voir Villiers p. 670) pour 19 M d’investissement."

$calcul.Range("A206").Value = "Armement moyen :"
$calcul.Range("D206").Formula = "=19000/198"

# -- comparison table: Villiers/Crowhurst investment & predicted captures ---
$calcul.Range("B209").Value = "Investissement Villiers / Crowhurst"
$calcul.Range("C209").Value = "Prises en utilisant le ratio de la guerre d’Amérpque"
$calcul.Range("H209").Value = "Valeur par prise : "
$calcul.Range("D209").Value = "Prises en utilisant Normann"

# -- Output sheet headers (year / French investment) need to exist before --
# -- "Part de captures privées" below, to reproduce the original authoring --
# -- order (and therefore shared-string table order) exactly.              --
$output.Range("A1").Value = "year"
$output.Range("C1").Value = "French investment"

# -- share of private captures out of the two totals (B136 vs B162) --------
$calcul.Range("D161").Value = "Part de captures privées"
$calcul.Range("D162").Formula = "=B136/(B136+B162)"
$calcul.Range("D162").NumberFormat = "#,##0.00"

$calcul.Range("E209").Value = "Ivestissement en utilisant Normann (que pour course privée)"
$output.Range("B1").Value = "French income"

$calcul.Range("I209").Formula = "=SUM(C210:C218)/(SUM(B169:B177)-B170-B171)"

$calcul.Range("A210").Value = 1793
$calcul.Range("B210").Formula = "=B195*$D$206"
$calcul.Range("C210").Formula = "=B210*$C$119"
$calcul.Range("D210").Formula = "=B169*$I$209"
$calcul.Range("E210").Formula = "=D210/$C$119*$D$162"

$calcul.Range("A211").Value = 1794
$calcul.Range("D211").Formula = "=B170*$I$209"
$calcul.Range("E211").Formula = "=D211/$C$119*$D$162"

$calcul.Range("A212").Value = 1795
$calcul.Range("D212").Formula = "=B171*$I$209"
$calcul.Range("E212").Formula = "=D212/$C$119*$D$162"

$calcul.Range("A213").Value = 1796
$calcul.Range("B213").Formula = "=B198*$D$206"
$calcul.Range("C213").Formula = "=B213*1.5"
$calcul.Range("D213").Formula = "=B172*$I$209"
$calcul.Range("E213").Formula = "=D213/$C$119*$D$162"

$calcul.Range("A214").Value = 1797
$calcul.Range("B214").Formula = "=B199*$D$206"
$calcul.Range("C214").Formula = "=B214*1.5"
$calcul.Range("D214").Formula = "=B173*$I$209"
$calcul.Range("E214").Formula = "=D214/$C$119*$D$162"

$calcul.Range("A215").Value = 1798
$calcul.Range("B215").Formula = "=B200*$D$206"
$calcul.Range("C215").Formula = "=B215*1.5"
$calcul.Range("D215").Formula = "=B174*$I$209"
$calcul.Range("E215").Formula = "=D215/$C$119*$D$162"

$calcul.Range("A216").Value = 1799
$calcul.Range("B216").Formula = "=B201*$D$206"
$calcul.Range("C216").Formula = "=B216*1.5"
$calcul.Range("D216").Formula = "=B175*$I$209"
$calcul.Range("E216").Formula = "=D216/$C$119*$D$162"

$calcul.Range("A217").Value = 1800
$calcul.Range("B217").Formula = "=B202*$D$206"
$calcul.Range("C217").Formula = "=B217*1.5"
$calcul.Range("D217").Formula = "=B176*$I$209"
$calcul.Range("E217").Formula = "=D217/$C$119*$D$162"

$calcul.Range("A218").Value = 1801
$calcul.Range("B218").Formula = "=B203*$D$206"
$calcul.Range("C218").Formula = "=B218*1.5"
$calcul.Range("D218").Formula = "=B177*$I$209"
$calcul.Range("E218").Formula = "=D218/$C$119*$D$162"

$calcul.Range("A219").Value = 1802
$calcul.Range("D219").Formula = "=B178*$I$209"
$calcul.Range("E219").Formula = "=D219/$C$119*$D$162"

$calcul.Range("A220").Value = 1803
$calcul.Range("D220").Formula = "=B179*$I$209"
$calcul.Range("E220").Formula = "=D220/$C$119*$D$162"

$calcul.Range("A221").Value = 1804
$calcul.Range("D221").Formula = "=B180*$I$209"
$calcul.Range("E221").Formula = "=D221/$C$119*$D$162"

$calcul.Range("A222").Value = 1805
$calcul.Range("D222").Formula = "=B181*$I$209"
$calcul.Range("E222").Formula = "=D222/$C$119*$D$162"

$calcul.Range("A223").Value = 1806
$calcul.Range("D223").Formula = "=B182*$I$209"
$calcul.Range("E223").Formula = "=D223/$C$119*$D$162"

$calcul.Range("A224").Value = 1807
$calcul.Range("D224").Formula = "=B183*$I$209"
$calcul.Range("E224").Formula = "=D224/$C$119*$D$162"

$calcul.Range("A225").Value = 1808
$calcul.Range("D225").Formula = "=B184*$I$209"
$calcul.Range("E225").Formula = "=D225/$C$119*$D$162"

$calcul.Range("A226").Value = 1809
$calcul.Range("D226").Formula = "=B185*$I$209"
$calcul.Range("E226").Formula = "=D226/$C$119*$D$162"

$calcul.Range("A227").Value = 1810
$calcul.Range("D227").Formula = "=B186*$I$209"
$calcul.Range("E227").Formula = "=D227/$C$119*$D$162"

$calcul.Range("A228").Value = 1811
$calcul.Range("D228").Formula = "=B187*$I$209"
$calcul.Range("E228").Formula = "=D228/$C$119*$D$162"

$calcul.Range("A229").Value = 1812
$calcul.Range("D229").Formula = "=B188*$I$209"
$calcul.Range("E229").Formula = "=D229/$C$119*$D$162"

$calcul.Range("A230").Value = 1813
$calcul.Range("D230").Formula = "=B189*$I$209"
$calcul.Range("E230").Formula = "=D230/$C$119*$D$162"

$calcul.Range("A231").Value = 1814
$calcul.Range("D231").Formula = "=B190*$I$209"
$calcul.Range("E231").Formula = "=D231/$C$119*$D$162"

# ============================================================================
# Output sheet
# ============================================================================

$output.Range("A1").Value = "year"
$output.Range("B1").Value = "French income"
$output.Range("C1").Value = "French investment"

$years = 1744..1814
foreach ($y in $years) {
    $r = $y - 1742
    $output.Range("A$r").Value = $y
}

$output.Range("B2").Formula = "=Calcul!B42"
$output.Range("C2").Formula = "=Calcul!B25"
$output.Range("B3").Formula = "=Calcul!B43"
$output.Range("C3").Formula = "=Calcul!B26"
$output.Range("B4").Formula = "=Calcul!B44"
$output.Range("C4").Formula = "=Calcul!B27"
$output.Range("B5").Formula = "=Calcul!B45"
$output.Range("C5").Formula = "=Calcul!B28"
$output.Range("B6").Formula = "=Calcul!B46"
$output.Range("C6").Formula = "=Calcul!B29"
$output.Range("B7").Value = $null
$output.Range("C7").Value = $null

$output.Range("B14").Formula = "=Calcul!B94+Calcul!D94"
$output.Range("C14").Formula = "=Calcul!B68"
$output.Range("B15").Formula = "=Calcul!B95+Calcul!D95"
$output.Range("C15").Formula = "=Calcul!B69"
$output.Range("B16").Formula = "=Calcul!B96+Calcul!D96"
$output.Range("C16").Formula = "=Calcul!B70"
$output.Range("B17").Formula = "=Calcul!B97+Calcul!D97"
$output.Range("C17").Formula = "=Calcul!B71"
$output.Range("B18").Formula = "=Calcul!B98+Calcul!D98"
$output.Range("C18").Formula = "=Calcul!B72"
$output.Range("B19").Formula = "=Calcul!B99+Calcul!D99"
$output.Range("C19").Formula = "=Calcul!B73"
$output.Range("B20").Formula = "=Calcul!B100+Calcul!D100"
$output.Range("C20").Formula = "=Calcul!B74"
$output.Range("B21").Value = $null
$output.Range("C21").Value = $null

$output.Range("B36").Formula = "=Calcul!B130+Calcul!B156"
$output.Range("C36").Formula = "=Calcul!B138"
$output.Range("B37").Formula = "=Calcul!B131+Calcul!B157"
$output.Range("C37").Formula = "=Calcul!B139"
$output.Range("B38").Formula = "=Calcul!B132+Calcul!B158"
$output.Range("C38").Formula = "=Calcul!B140"
$output.Range("B39").Formula = "=Calcul!B133+Calcul!B159"
$output.Range("C39").Formula = "=Calcul!B141"
$output.Range("B40").Formula = "=Calcul!B134+Calcul!B160"
$output.Range("C40").Formula = "=Calcul!B142"
$output.Range("B41").Formula = "=Calcul!B135+Calcul!B161"
$output.Range("C41").Formula = "=Calcul!B143"
$output.Range("B42").Value = $null

$output.Range("B51").Formula = "=Calcul!D210"
$output.Range("C51").Formula = "=Calcul!E210"
$output.Range("B52").Formula = "=Calcul!D211"
$output.Range("C52").Formula = "=Calcul!E211"
$output.Range("B53").Formula = "=Calcul!D212"
$output.Range("C53").Formula = "=Calcul!E212"
$output.Range("B54").Formula = "=Calcul!D213"
$output.Range("C54").Formula = "=Calcul!E213"
$output.Range("B55").Formula = "=Calcul!D214"
$output.Range("C55").Formula = "=Calcul!E214"
$output.Range("B56").Formula = "=Calcul!D215"
$output.Range("C56").Formula = "=Calcul!E215"
$output.Range("B57").Formula = "=Calcul!D216"
$output.Range("C57").Formula = "=Calcul!E216"
$output.Range("B58").Formula = "=Calcul!D217"
$output.Range("C58").Formula = "=Calcul!E217"
$output.Range("B59").Formula = "=Calcul!D218"
$output.Range("C59").Formula = "=Calcul!E218"
$output.Range("B60").Value = $null
$output.Range("C60").Value = $null

$output.Range("B61").Formula = "=Calcul!D220"
$output.Range("C61").Formula = "=Calcul!E220"
$output.Range("B62").Formula = "=Calcul!D221"
$output.Range("C62").Formula = "=Calcul!E221"
$output.Range("B63").Formula = "=Calcul!D222"
$output.Range("C63").Formula = "=Calcul!E222"
$output.Range("B64").Formula = "=Calcul!D223"
$output.Range("C64").Formula = "=Calcul!E223"
$output.Range("B65").Formula = "=Calcul!D224"
$output.Range("C65").Formula = "=Calcul!E224"
$output.Range("B66").Formula = "=Calcul!D225"
$output.Range("C66").Formula = "=Calcul!E225"
$output.Range("B67").Formula = "=Calcul!D226"
$output.Range("C67").Formula = "=Calcul!E226"
$output.Range("B68").Formula = "=Calcul!D227"
$output.Range("C68").Formula = "=Calcul!E227"
$output.Range("B69").Formula = "=Calcul!D228"
$output.Range("C69").Formula = "=Calcul!E228"
$output.Range("B70").Formula = "=Calcul!D229"
$output.Range("C70").Formula = "=Calcul!E229"
$output.Range("B71").Formula = "=Calcul!D230"
$output.Range("C71").Formula = "=Calcul!E230"
$output.Range("B72").Formula = "=Calcul!D231"
$output.Range("C72").Formula = "=Calcul!E231"

# -- selection / active sheet, matching the final screen state --------------
$calcul.Range("C94").Select()
$output.Activate()
$output.Range("C61").Select()
